# Apply updated crypto price/volume data per the commit diff.
# Cells in column D whose new text looks like a plain number are
# prefixed with a leading apostrophe so Excel keeps them as text
# (matching the original inlineStr/text cell type) instead of
# auto-converting them to a numeric value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.249.09"
$ws.Range("E2").Value = "  +2.01%  "
$ws.Range("D3").Value = "3.926.53"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'487.10"
$ws.Range("E5").Value = "  +3.26%  "
$ws.Range("D6").Value = "'148.77"
$ws.Range("E6").Value = "  +3.40%  "
$ws.Range("E7").Value = "  +1.32%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "'0.734"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").Value = "'0.169"
$ws.Range("E10").Value = "  +3.75%  "
$ws.Range("D11").Value = "'0.0000352"
$ws.Range("E11").Value = "  +5.39%  "
$ws.Range("D12").Value = "'43.09"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").Value = "'10.68"
$ws.Range("E13").Value = "  +3.29%  "
$ws.Range("D14").Value = "4.547.27"
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("D15").Value = "'14.79"
$ws.Range("E15").Value = "  -0.75%  "
$ws.Range("D16").Value = "3.911.78"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").Value = "'20.03"
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("D19").Value = "'1.14"
$ws.Range("E19").Value = "  -1.82%  "
$ws.Range("D20").Value = "68.355.45"
$ws.Range("E20").Value = "  +1.85%  "
$ws.Range("D21").Value = "'440.70"
$ws.Range("E21").Value = "  +2.35%  "
$ws.Range("B22").Value = "ImmutableX"
$ws.Range("C22").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D22").Value = "'3.45"
$ws.Range("E22").Value = "  +2.48%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Value = "'14.96"
$ws.Range("E23").Value = "  +2.64%  "
$ws.Range("D24").Value = "'88.62"
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("D25").Value = "'11.52"
$ws.Range("E25").Value = "  +15.36%  "
$ws.Range("D26").Value = "'11.18"
$ws.Range("E26").Value = "  +16.18%  "
$ws.Range("D27").Value = "'3.64"
$ws.Range("E27").Value = "  +2.48%  "
$ws.Range("D28").Value = "'38.74"
$ws.Range("E28").Value = "  +1.21%  "
$ws.Range("D29").Value = "'5.87"
$ws.Range("E29").Value = "  +1.70%  "
$ws.Range("D30").Value = "'718.96"
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("E31").Value = "  -1.22%  "
$ws.Range("E32").Value = "  -0.39%  "
$ws.Range("D33").Value = "'2.89"
$ws.Range("E33").Value = "  +3.10%  "
$ws.Range("D34").Value = "0.0₃0913"
$ws.Range("E34").Value = "  +15.40%  "
$ws.Range("D35").Value = "'41.89"
$ws.Range("E35").Value = "  -3.29%  "
$ws.Range("D36").Value = "'6.16"
$ws.Range("E36").Value = "  +13.82%  "
$ws.Range("D37").Value = "'59.41"
$ws.Range("E37").Value = "  +3.83%  "
$ws.Range("E38").Value = "  -3.87%  "
$ws.Range("D39").Value = "'0.398"
$ws.Range("E39").Value = "  +18.04%  "
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").Value = "'2.98"
$ws.Range("E41").Value = "  +18.17%  "
$ws.Range("E42").Value = "  +1.69%  "
$ws.Range("D43").Value = "'3.19"
$ws.Range("E43").Value = "  +3.63%  "
$ws.Range("D44").Value = "'2.94"
$ws.Range("E44").Value = "  +5.65%  "
$ws.Range("E45").Value = "  +1.58%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "0.0₆0362"
$ws.Range("E47").Value = "  +47.17%  "
$ws.Range("D48").Value = "'3.42"
$ws.Range("E48").Value = "  +1.04%  "
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("D50").Value = "'145.52"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("D51").Value = "'3.14"
$ws.Range("E51").Value = "  +0.18%  "
